# Insert a new data row at row 363 (shifts existing rows 363:470 down to 364:471)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("363:363").Insert()

# Populate the newly inserted row 363 with the new weekly price observation
$ws.Cells.Item(363, 1).Value = 4
$ws.Cells.Item(363, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(363, 3).Value = "Los Lagos"
$ws.Cells.Item(363, 4).Value = 45093
$ws.Cells.Item(363, 5).Value = 10
$ws.Cells.Item(363, 6).Value = 100112037
$ws.Cells.Item(363, 7).Value = "Cebollín"
$ws.Cells.Item(363, 8).Value = "Sin especificar"
$ws.Cells.Item(363, 9).Value = "Primera"
$ws.Cells.Item(363, 10).Value = 180
$ws.Cells.Item(363, 11).Value = 6500
$ws.Cells.Item(363, 12).Value = 6500
$ws.Cells.Item(363, 13).Value = 6500
$ws.Cells.Item(363, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(363, 15).Value = "Región Metropolitana"
$ws.Cells.Item(363, 16).Value = 181
$ws.Cells.Item(363, 17).Value = 36
$ws.Cells.Item(363, 18).Value = "Hortaliza"
